$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Thu Jan 25 17:18:16 EST 2024"
$ws.Range("B3").Value = "Thu Jan 25 17:18:26 EST 2024"
$ws.Range("B4").Value = "Thu Jan 25 17:18:36 EST 2024"
$ws.Range("B5").Value = "Thu Jan 25 17:18:46 EST 2024"
$ws.Range("B6").Value = "Thu Jan 25 17:18:57 EST 2024"
$ws.Range("B7").Value = "Thu Jan 25 17:19:07 EST 2024"
$ws.Range("B8").Value = "Thu Jan 25 17:19:17 EST 2024"
$ws.Range("B9").Value = "Thu Jan 25 17:19:27 EST 2024"
$ws.Range("B10").Value = "Thu Jan 25 17:19:37 EST 2024"
$ws.Range("B13").Value = "Thu Jan 25 17:19:47 EST 2024"
$ws.Range("B14").Value = "Thu Jan 25 17:19:57 EST 2024"
$ws.Range("B15").Value = "Thu Jan 25 17:20:07 EST 2024"
$ws.Range("B16").Value = "Thu Jan 25 17:20:18 EST 2024"
$ws.Range("B17").Value = "Thu Jan 25 17:20:27 EST 2024"
$ws.Range("B18").Value = "Thu Jan 25 17:20:38 EST 2024"
$ws.Range("B19").Value = "Thu Jan 25 17:20:48 EST 2024"
$ws.Range("B20").Value = "Thu Jan 25 17:20:58 EST 2024"
$ws.Range("B21").Value = "Thu Jan 25 17:21:08 EST 2024"
$ws.Range("B22").Value = "Thu Jan 25 17:21:18 EST 2024"
$ws.Range("B23").Value = "Thu Jan 25 17:21:28 EST 2024"
$ws.Range("B24").Value = "Thu Jan 25 17:21:39 EST 2024"
$ws.Range("B25").Value = "Thu Jan 25 17:21:49 EST 2024"
$ws.Range("B26").Value = "Thu Jan 25 17:22:00 EST 2024"
$ws.Range("B27").Value = "Thu Jan 25 17:22:10 EST 2024"
$ws.Range("B28").Value = "Thu Jan 25 17:22:20 EST 2024"
$ws.Range("B29").Value = "Thu Jan 25 17:22:30 EST 2024"
$ws.Range("B30").Value = "Thu Jan 25 17:22:41 EST 2024"
